# Grump metadata. Longhurst provinces and Seasons added.
# Add three new variable-metadata rows (Longhurst_Long, Longhurst_Short, Season)
# to the "vars_meta_data" sheet, right below the existing rows of data.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("vars_meta_data")
$ws.Activate()

# Columns: A=var_short_name B=var_long_name C=var_sensor D=var_unit
#          E=var_spatial_res F=var_temporal_res G=var_discipline H=visualize

$rows = @(
  @{ Row = 41; Short = "Longhurst_Long";  Long = "Longhurst province sample was taken in." },
  @{ Row = 42; Short = "Longhurst_Short"; Long = "Longhurst province sample was taken in, shortened code." },
  @{ Row = 43; Short = "Season";          Long = "Season sample was taken in." }
)

foreach ($r in $rows) {
  $i = $r.Row
  $ws.Cells.Item($i, 1).Value = $r.Short
  $ws.Cells.Item($i, 2).Value = $r.Long
  $ws.Cells.Item($i, 3).Value = "NA"
  $ws.Cells.Item($i, 4).Value = "NA"
  $ws.Cells.Item($i, 5).Value = "Irregular"
  $ws.Cells.Item($i, 6).Value = "Irregular"
  $ws.Cells.Item($i, 7).Value = "Biology"
  $ws.Cells.Item($i, 8).Value = 1
}

# Scroll the view up so row 17 is at the top, then select the newly added rows,
# matching where the author's cursor ended up after data entry.
$win = $excel.ActiveWindow
$win.ScrollRow = 17
$win.ScrollColumn = 1
$ws.Range("A41:H43").Select()
